# Auto-generated Excel COM-interop script applying the Tonberry_Profits
# leve-profit recalculation update across all eight job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$edits_ALC = @(
    @{Cell="H18"; Action="set"; Value=16126.286},
    @{Cell="I18"; Action="set"; Value=7000},
    @{Cell="J18"; Action="set"; Value=17647.334},
    @{Cell="K18"; Action="set"; Value=7000},
    @{Cell="L18"; Action="set"; Value=17647.334},
    @{Cell="M18"; Action="set"; Value=-6716},
    @{Cell="N18"; Action="set"; Value=-18215.334},
    @{Cell="H21"; Action="set"; Value=59679.332},
    @{Cell="I21"; Action="set"; Value=59679.332},
    @{Cell="K21"; Action="set"; Value=59679.332},
    @{Cell="M21"; Action="set"; Value=-59211.332},
    @{Cell="H23"; Action="set"; Value=59679.332},
    @{Cell="I23"; Action="set"; Value=59679.332},
    @{Cell="K23"; Action="set"; Value=59679.332},
    @{Cell="M23"; Action="set"; Value=-59445.332},
    @{Cell="H53"; Action="set"; Value=7949.615},
    @{Cell="I53"; Action="set"; Value=11337.223},
    @{Cell="K53"; Action="set"; Value=11337.223},
    @{Cell="M53"; Action="set"; Value=-10700.223},
    @{Cell="H100"; Action="set"; Value=2684.1667},
    @{Cell="I100"; Action="set"; Value=2821},
    @{Cell="K100"; Action="set"; Value=2821},
    @{Cell="M100"; Action="set"; Value=-2280},
    @{Cell="H129"; Action="set"; Value=850.0540999999999},
    @{Cell="I129"; Action="set"; Value=631.1667},
    @{Cell="J129"; Action="set"; Value=869.3677},
    @{Cell="K129"; Action="set"; Value=1893.5001},
    @{Cell="L129"; Action="set"; Value=2608.1031},
    @{Cell="M129"; Action="set"; Value=3106.4999},
    @{Cell="N129"; Action="set"; Value=-12608.1031},
    @{Cell="H132"; Action="set"; Value=1074.25},
    @{Cell="I132"; Action="set"; Value=1016.64},
    @{Cell="J132"; Action="set"; Value=1280},
    @{Cell="K132"; Action="set"; Value=3049.92},
    @{Cell="L132"; Action="set"; Value=3840},
    @{Cell="M132"; Action="set"; Value=-519.9200000000001},
    @{Cell="N132"; Action="set"; Value=-8900},
    @{Cell="H135"; Action="set"; Value=504.6},
    @{Cell="I135"; Action="set"; Value=504.6},
    @{Cell="J135"; Action="set"; Value=0},
    @{Cell="K135"; Action="set"; Value=4541.400000000001},
    @{Cell="L135"; Action="set"; Value=0},
    @{Cell="M135"; Action="set"; Value=-2006.400000000001},
    @{Cell="N135"; Action="delete"},
    @{Cell="H137"; Action="set"; Value=2019.6875},
    @{Cell="J137"; Action="set"; Value=2294.125},
    @{Cell="L137"; Action="set"; Value=6882.375},
    @{Cell="N137"; Action="set"; Value=-11982.375},
    @{Cell="H138"; Action="set"; Value=1812.1587},
    @{Cell="I138"; Action="set"; Value=1358.4231},
    @{Cell="K138"; Action="set"; Value=4075.2693},
    @{Cell="M138"; Action="set"; Value=1064.7307},
    @{Cell="H141"; Action="set"; Value=3941.9375},
    @{Cell="I141"; Action="set"; Value=3152.7},
    @{Cell="J141"; Action="set"; Value=5257.3335},
    @{Cell="K141"; Action="set"; Value=9458.099999999999},
    @{Cell="L141"; Action="set"; Value=15772.0005},
    @{Cell="M141"; Action="set"; Value=-4278.099999999999},
    @{Cell="N141"; Action="set"; Value=-26132.0005}
)
$ws = $wb.Worksheets.Item("ALC")
foreach ($e in $edits_ALC) {
    if ($e.Action -eq "delete") {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value = $e.Value
    }
}

# --- ARM sheet ---
$edits_ARM = @(
    @{Cell="H32"; Action="set"; Value=3211.739},
    @{Cell="I32"; Action="set"; Value=1900.7119},
    @{Cell="K32"; Action="set"; Value=1900.7119},
    @{Cell="M32"; Action="set"; Value=-1613.7119},
    @{Cell="H122"; Action="set"; Value=2206},
    @{Cell="I122"; Action="set"; Value=2206},
    @{Cell="K122"; Action="set"; Value=6618},
    @{Cell="M122"; Action="set"; Value=-4168},
    @{Cell="H132"; Action="set"; Value=1739.5518},
    @{Cell="I132"; Action="set"; Value=1184.7826},
    @{Cell="K132"; Action="set"; Value=3554.3478},
    @{Cell="M132"; Action="set"; Value=-1024.3478}
)
$ws = $wb.Worksheets.Item("ARM")
foreach ($e in $edits_ARM) {
    if ($e.Action -eq "delete") {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value = $e.Value
    }
}

# --- BSM sheet ---
$edits_BSM = @(
    @{Cell="H62"; Action="set"; Value=10000},
    @{Cell="J62"; Action="set"; Value=0},
    @{Cell="L62"; Action="set"; Value=0},
    @{Cell="N62"; Action="delete"},
    @{Cell="H65"; Action="set"; Value=10000},
    @{Cell="J65"; Action="set"; Value=0},
    @{Cell="L65"; Action="set"; Value=0},
    @{Cell="N65"; Action="delete"},
    @{Cell="H86"; Action="set"; Value=89640.78},
    @{Cell="J86"; Action="set"; Value=144711},
    @{Cell="L86"; Action="set"; Value=144711},
    @{Cell="N86"; Action="set"; Value=-146957},
    @{Cell="H89"; Action="set"; Value=89640.78},
    @{Cell="J89"; Action="set"; Value=144711},
    @{Cell="L89"; Action="set"; Value=723555},
    @{Cell="N89"; Action="set"; Value=-734787},
    @{Cell="H134"; Action="set"; Value=4437.154},
    @{Cell="I134"; Action="set"; Value=4437.154},
    @{Cell="K134"; Action="set"; Value=13311.462},
    @{Cell="M134"; Action="set"; Value=-10776.462}
)
$ws = $wb.Worksheets.Item("BSM")
foreach ($e in $edits_BSM) {
    if ($e.Action -eq "delete") {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value = $e.Value
    }
}

# --- CRP sheet ---
$edits_CRP = @(
    @{Cell="H10"; Action="set"; Value=1342.8572},
    @{Cell="I10"; Action="set"; Value=800},
    @{Cell="J10"; Action="set"; Value=2700},
    @{Cell="K10"; Action="set"; Value=800},
    @{Cell="L10"; Action="set"; Value=2700},
    @{Cell="M10"; Action="set"; Value=-661},
    @{Cell="N10"; Action="set"; Value=-2978},
    @{Cell="H19"; Action="set"; Value=1116.25},
    @{Cell="I19"; Action="set"; Value=1153.3334},
    @{Cell="J19"; Action="set"; Value=1005},
    @{Cell="K19"; Action="set"; Value=1153.3334},
    @{Cell="L19"; Action="set"; Value=1005},
    @{Cell="M19"; Action="set"; Value=-983.3334},
    @{Cell="N19"; Action="set"; Value=-1345},
    @{Cell="H24"; Action="set"; Value=1116.25},
    @{Cell="I24"; Action="set"; Value=1153.3334},
    @{Cell="J24"; Action="set"; Value=1005},
    @{Cell="K24"; Action="set"; Value=1153.3334},
    @{Cell="L24"; Action="set"; Value=1005},
    @{Cell="M24"; Action="set"; Value=-983.3334},
    @{Cell="N24"; Action="set"; Value=-1345},
    @{Cell="H31"; Action="set"; Value=2202.7666},
    @{Cell="I31"; Action="set"; Value=1828.6},
    @{Cell="J31"; Action="set"; Value=2576.9333},
    @{Cell="K31"; Action="set"; Value=1828.6},
    @{Cell="L31"; Action="set"; Value=2576.9333},
    @{Cell="M31"; Action="set"; Value=-1533.6},
    @{Cell="N31"; Action="set"; Value=-3166.9333},
    @{Cell="H34"; Action="set"; Value=2202.7666},
    @{Cell="I34"; Action="set"; Value=1828.6},
    @{Cell="J34"; Action="set"; Value=2576.9333},
    @{Cell="K34"; Action="set"; Value=1828.6},
    @{Cell="L34"; Action="set"; Value=2576.9333},
    @{Cell="M34"; Action="set"; Value=-1626.6},
    @{Cell="N34"; Action="set"; Value=-2980.9333},
    @{Cell="H58"; Action="set"; Value=1978138},
    @{Cell="I58"; Action="set"; Value=3106867.8},
    @{Cell="K58"; Action="set"; Value=3106867.8},
    @{Cell="M58"; Action="set"; Value=-3106664.8},
    @{Cell="H105"; Action="set"; Value=2302.5},
    @{Cell="I105"; Action="set"; Value=2403.3333},
    @{Cell="J105"; Action="set"; Value=2000},
    @{Cell="K105"; Action="set"; Value=2403.3333},
    @{Cell="L105"; Action="set"; Value=2000},
    @{Cell="M105"; Action="set"; Value=-656.3332999999998},
    @{Cell="N105"; Action="set"; Value=-5494},
    @{Cell="H122"; Action="set"; Value=10014},
    @{Cell="I122"; Action="set"; Value=0},
    @{Cell="K122"; Action="set"; Value=0},
    @{Cell="M122"; Action="delete"},
    @{Cell="H132"; Action="set"; Value=2198.147},
    @{Cell="I132"; Action="set"; Value=1266.25},
    @{Cell="J132"; Action="set"; Value=3529.4285},
    @{Cell="K132"; Action="set"; Value=3798.75},
    @{Cell="L132"; Action="set"; Value=10588.2855},
    @{Cell="M132"; Action="set"; Value=-1268.75},
    @{Cell="N132"; Action="set"; Value=-15648.2855},
    @{Cell="H134"; Action="set"; Value=1249.8846},
    @{Cell="I134"; Action="set"; Value=1249.8846},
    @{Cell="K134"; Action="set"; Value=3749.6538},
    @{Cell="M134"; Action="set"; Value=-1214.6538},
    @{Cell="H136"; Action="set"; Value=1978138},
    @{Cell="I136"; Action="set"; Value=3106867.8},
    @{Cell="K136"; Action="set"; Value=9320603.399999999},
    @{Cell="M136"; Action="set"; Value=-9318053.399999999}
)
$ws = $wb.Worksheets.Item("CRP")
foreach ($e in $edits_CRP) {
    if ($e.Action -eq "delete") {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value = $e.Value
    }
}

# --- CUL sheet ---
$edits_CUL = @(
    @{Cell="H11"; Action="set"; Value=566.3333},
    @{Cell="J11"; Action="set"; Value=566.3333},
    @{Cell="L11"; Action="set"; Value=1698.9999},
    @{Cell="N11"; Action="set"; Value=-1978.9999},
    @{Cell="H114"; Action="set"; Value=2175.8333},
    @{Cell="I114"; Action="set"; Value=585.3333},
    @{Cell="J114"; Action="set"; Value=3766.3333},
    @{Cell="K114"; Action="set"; Value=1755.9999},
    @{Cell="L114"; Action="set"; Value=11298.9999},
    @{Cell="M114"; Action="set"; Value=1498.0001},
    @{Cell="N114"; Action="set"; Value=-17806.9999},
    @{Cell="H131"; Action="set"; Value=780.98},
    @{Cell="J131"; Action="set"; Value=800.4457},
    @{Cell="L131"; Action="set"; Value=2401.3371},
    @{Cell="N131"; Action="set"; Value=-12481.3371},
    @{Cell="H137"; Action="set"; Value=2972.5},
    @{Cell="J137"; Action="set"; Value=3493.3333},
    @{Cell="L137"; Action="set"; Value=10479.9999},
    @{Cell="N137"; Action="set"; Value=-20679.9999},
    @{Cell="H139"; Action="set"; Value=10895.818},
    @{Cell="I139"; Action="set"; Value=12650.444},
    @{Cell="J139"; Action="set"; Value=3000},
    @{Cell="K139"; Action="set"; Value=37951.33199999999},
    @{Cell="L139"; Action="set"; Value=9000},
    @{Cell="M139"; Action="set"; Value=-32811.33199999999},
    @{Cell="N139"; Action="set"; Value=-19280}
)
$ws = $wb.Worksheets.Item("CUL")
foreach ($e in $edits_CUL) {
    if ($e.Action -eq "delete") {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value = $e.Value
    }
}

# --- GSM sheet ---
$edits_GSM = @(
    @{Cell="H122"; Action="set"; Value=2314.5386},
    @{Cell="I122"; Action="set"; Value=1832},
    @{Cell="J122"; Action="set"; Value=2728.1428},
    @{Cell="K122"; Action="set"; Value=5496},
    @{Cell="L122"; Action="set"; Value=8184.428400000001},
    @{Cell="M122"; Action="set"; Value=-3046},
    @{Cell="N122"; Action="set"; Value=-13084.4284},
    @{Cell="H132"; Action="set"; Value=3207245.8},
    @{Cell="I132"; Action="set"; Value=3847806},
    @{Cell="J132"; Action="set"; Value=4443.5},
    @{Cell="K132"; Action="set"; Value=11543418},
    @{Cell="L132"; Action="set"; Value=13330.5},
    @{Cell="M132"; Action="set"; Value=-11540888},
    @{Cell="N132"; Action="set"; Value=-18390.5}
)
$ws = $wb.Worksheets.Item("GSM")
foreach ($e in $edits_GSM) {
    if ($e.Action -eq "delete") {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value = $e.Value
    }
}

# --- LTW sheet ---
$edits_LTW = @(
    @{Cell="H16"; Action="set"; Value=3472.4375},
    @{Cell="I16"; Action="set"; Value=4245.125},
    @{Cell="J16"; Action="set"; Value=2699.75},
    @{Cell="K16"; Action="set"; Value=4245.125},
    @{Cell="L16"; Action="set"; Value=2699.75},
    @{Cell="M16"; Action="set"; Value=-4075.125},
    @{Cell="N16"; Action="set"; Value=-3039.75},
    @{Cell="H40"; Action="set"; Value=2471.1428},
    @{Cell="I40"; Action="set"; Value=2459.8},
    @{Cell="K40"; Action="set"; Value=2459.8},
    @{Cell="M40"; Action="set"; Value=-2323.8},
    @{Cell="H132"; Action="set"; Value=3058},
    @{Cell="I132"; Action="set"; Value=2303.3125},
    @{Cell="K132"; Action="set"; Value=6909.9375},
    @{Cell="M132"; Action="set"; Value=-4379.9375},
    @{Cell="H136"; Action="set"; Value=3731.4814},
    @{Cell="I136"; Action="set"; Value=3102},
    @{Cell="K136"; Action="set"; Value=9306},
    @{Cell="M136"; Action="set"; Value=-6756}
)
$ws = $wb.Worksheets.Item("LTW")
foreach ($e in $edits_LTW) {
    if ($e.Action -eq "delete") {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value = $e.Value
    }
}

# --- WVR sheet ---
$edits_WVR = @(
    @{Cell="H108"; Action="set"; Value=67999.5},
    @{Cell="J108"; Action="set"; Value=67999.5},
    @{Cell="L108"; Action="set"; Value=67999.5},
    @{Cell="N108"; Action="set"; Value=-75679.5},
    @{Cell="H132"; Action="set"; Value=1477.1143},
    @{Cell="I132"; Action="set"; Value=970.2308},
    @{Cell="J132"; Action="set"; Value=2941.4443},
    @{Cell="K132"; Action="set"; Value=2910.6924},
    @{Cell="L132"; Action="set"; Value=8824.332900000001},
    @{Cell="M132"; Action="set"; Value=-380.6923999999999},
    @{Cell="N132"; Action="set"; Value=-13884.3329},
    @{Cell="H136"; Action="set"; Value=23150576},
    @{Cell="I136"; Action="set"; Value=32681890},
    @{Cell="K136"; Action="set"; Value=98045670},
    @{Cell="M136"; Action="set"; Value=-98043120}
)
$ws = $wb.Worksheets.Item("WVR")
foreach ($e in $edits_WVR) {
    if ($e.Action -eq "delete") {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value = $e.Value
    }
}

